$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 29
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05).
for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
